# Refresh the coin Price (D) and Volume(1h) (E) columns with the latest
# snapshot values, as produced by the scheduled cryptos-list update workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell address -> new text value.
$updates = [ordered]@{
    "D2" = "27.513.73"
    "E2" = "  -0.60%  "
    "D3" = "1.621.41"
    "E3" = "  -1.31%  "
    "E4" = "  -0.04%  "
    "D5" = "211.51"
    "E5" = "  -0.73%  "
    "E6" = "  -0.94%  "
    "E7" = "  -0.04%  "
    "D8" = "23.02"
    "E8" = "  -1.22%  "
    "E9" = "  +1.67%  "
    "E10" = "  -0.13%  "
    "D11" = "0.0882"
    "E11" = "  -1.36%  "
    "D12" = "1.850.93"
    "D13" = "1.608.76"
    "E13" = "  -2.11%  "
    "E14" = "  -0.21%  "
    "E15" = "  -2.34%  "
    "D16" = "65.28"
    "E16" = "  +0.89%  "
    "D17" = "27.485.08"
    "D18" = "229.83"
    "E18" = "  -0.60%  "
    "E19" = "  -1.00%  "
    "E20" = "  -1.82%  "
    "E21" = "  +0.02%  "
    "D22" = "10.37"
    "E22" = "  +3.33%  "
    "E23" = "  +0.97%  "
    "E24" = "  +7.61%  "
    "D25" = "149.84"
    "E25" = "  -0.07%  "
    "E26" = "  -1.08%  "
    "E27" = "  -1.08%  "
    "E28" = "  +0.00%  "
    "E29" = "  -0.96%  "
    "E30" = "  -1.08%  "
    "E31" = "  -0.79%  "
    "E32" = "  -1.15%  "
    "D33" = "1.466.15"
    "E33" = "  +0.80%  "
    "E34" = "  -2.80%  "
    "E35" = "  -2.18%  "
    "E36" = "  -0.42%  "
    "D37" = "0.946"
    "E37" = "  +4.22%  "
    "D38" = "0.0167"
    "E38" = "  -0.28%  "
    "E39" = "  -1.39%  "
    "D40" = "0.553"
    "E40" = "  -3.03%  "
    "E41" = "  -0.02%  "
    "E42" = "  -2.37%  "
    "D43" = "67.29"
    "E43" = "  -4.97%  "
    "D44" = "2.49"
    "E44" = "  +0.41%  "
    "E45" = "  -2.23%  "
    "D46" = "5.31"
    "E46" = "  -5.88%  "
    "E47" = "  +1.76%  "
    "D48" = "1.760.87"
    "E48" = "  -1.50%  "
    "D49" = "87.27"
    "E49" = "  +1.39%  "
    "D50" = "0.0₆0105"
    "E50" = "  -1.55%  "
    "D51" = "0.0996"
    "E51" = "  +0.46%  "
}

# Cells whose new value would otherwise be auto-parsed by Excel as a plain
# number (e.g. "211.51") need to be forced to Text so they stay strings,
# just like every other cell in these columns (originally inline strings).
$forceTextCells = @(
    "D5", "D8", "D11", "D16", "D18", "D22", "D25", "D37", "D38", "D40", "D43", "D44", "D46", "D49", "D51"
)

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    if ($forceTextCells -contains $addr) {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
